# Applies the benchmark-table refresh described by the commit:
# "Fixed README.md stats and docx preparation for all DaCapo - JDK 21 - Z GC tests"
#
# The document is a single one-column table; each row holds one reported
# statistic. We update the affected rows' text in place (preserving each
# row's existing run/cell formatting, which is identical across every row),
# and collapse the three multi-run "tab separated" summary rows down to
# their single trailing value - matching the target OOXML exactly without
# needing to insert/delete rows (row count is unchanged end to end).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-RowText($rowIndex, $newText) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $newText
}

Set-RowText 1 "0M"
Set-RowText 2 "0M"
Set-RowText 3 "0M"
Set-RowText 4 "56"
Set-RowText 5 "0.00003"
Set-RowText 6 "0.00011"
Set-RowText 8 "0.00002"
Set-RowText 9 "0.00008"
Set-RowText 10 "0.00009"
Set-RowText 11 "0.00010"
Set-RowText 12 "0.00470"

Set-RowText 44 "100"
Set-RowText 45 "0"
Set-RowText 46 "229"
